# Update NATMI ligand-receptor edge weights with re-computed TPM values.
# Columns G,H,I,J (ligand avg/total expr + specificity) and
# M,N,O,P (receptor avg/total expr + specificity) and
# Q,R,S,T (derived edge weights/specificities) are refreshed per row;
# identifier columns A-F, K, L are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 78.90112033333334
$ws.Range("H2").Value = 236.703361
$ws.Range("I2").Value = 0.3371779636489425
$ws.Range("J2").Value = 0.3371779636489425
$ws.Range("M2").Value = 5.375839
$ws.Range("N2").Value = 16.127517
$ws.Range("O2").Value = 0.2354568587499626
$ws.Range("P2").Value = 0.2354568587499626
$ws.Range("Q2").Value = 424.1597198316264
$ws.Range("R2").Value = 3817.437478484637
$ws.Range("S2").Value = 0.07939086416048907
$ws.Range("T2").Value = 0.07939086416048909
# Row 3
$ws.Range("G3").Value = 78.90112033333334
$ws.Range("H3").Value = 236.703361
$ws.Range("I3").Value = 0.3371779636489425
$ws.Range("J3").Value = 0.3371779636489425
$ws.Range("O3").Value = 0.007131134316291014
$ws.Range("P3").Value = 0.007131134316291014
$ws.Range("Q3").Value = 12.84625960669822
$ws.Range("R3").Value = 115.616336460284
$ws.Range("S3").Value = 0.002404461347274098
$ws.Range("T3").Value = 0.002404461347274098
# Row 4
$ws.Range("G4").Value = 78.90112033333334
$ws.Range("H4").Value = 236.703361
$ws.Range("I4").Value = 0.3371779636489425
$ws.Range("J4").Value = 0.3371779636489425
$ws.Range("M4").Value = 9.994147
$ws.Range("N4").Value = 29.982441
$ws.Range("O4").Value = 0.4377345486919088
$ws.Range("P4").Value = 0.4377345486919088
$ws.Range("Q4").Value = 788.5493950760224
$ws.Range("R4").Value = 7096.9445556842
$ws.Range("S4").Value = 0.1475944437467267
$ws.Range("T4").Value = 0.1475944437467267
# Row 5
$ws.Range("G5").Value = 78.90112033333334
$ws.Range("H5").Value = 236.703361
$ws.Range("I5").Value = 0.3371779636489425
$ws.Range("J5").Value = 0.3371779636489425
$ws.Range("M5").Value = 0.7761303333333333
$ws.Range("N5").Value = 2.328391
$ws.Range("O5").Value = 0.03399380269149206
$ws.Range("P5").Value = 0.03399380269149207
$ws.Range("Q5").Value = 61.23755282468344
$ws.Range("R5").Value = 551.1379754221509
$ws.Range("S5").Value = 0.01146196116820124
$ws.Range("T5").Value = 0.01146196116820124
# Row 6
$ws.Range("G6").Value = 78.90112033333334
$ws.Range("H6").Value = 236.703361
$ws.Range("I6").Value = 0.3371779636489425
$ws.Range("J6").Value = 0.3371779636489425
$ws.Range("M6").Value = 6.522593333333333
$ws.Range("N6").Value = 19.56778
$ws.Range("O6").Value = 0.2856836555503455
$ws.Range("P6").Value = 0.2856836555503455
$ws.Range("Q6").Value = 514.6399214787311
$ws.Range("R6").Value = 4631.75929330858
$ws.Range("S6").Value = 0.09632623322625143
$ws.Range("T6").Value = 0.09632623322625143
# Row 7
$ws.Range("I7").Value = 0.1683260544097508
$ws.Range("J7").Value = 0.1683260544097508
$ws.Range("M7").Value = 5.375839
$ws.Range("N7").Value = 16.127517
$ws.Range("O7").Value = 0.2354568587499626
$ws.Range("P7").Value = 0.2354568587499626
$ws.Range("Q7").Value = 211.749105149526
$ws.Range("R7").Value = 1905.741946345734
$ws.Range("S7").Value = 0.0396335240170952
$ws.Range("T7").Value = 0.03963352401709521
# Row 8
$ws.Range("I8").Value = 0.1683260544097508
$ws.Range("J8").Value = 0.1683260544097508
$ws.Range("O8").Value = 0.007131134316291014
$ws.Range("P8").Value = 0.007131134316291014
$ws.Range("S8").Value = 0.001200355702927242
$ws.Range("T8").Value = 0.001200355702927242
# Row 9
$ws.Range("I9").Value = 0.1683260544097508
$ws.Range("J9").Value = 0.1683260544097508
$ws.Range("M9").Value = 9.994147
$ws.Range("N9").Value = 29.982441
$ws.Range("O9").Value = 0.4377345486919088
$ws.Range("P9").Value = 0.4377345486919088
$ws.Range("Q9").Value = 393.659795983998
$ws.Range("R9").Value = 3542.938163855982
$ws.Range("S9").Value = 0.07368212946014194
$ws.Range("T9").Value = 0.07368212946014194
# Row 10
$ws.Range("I10").Value = 0.1683260544097508
$ws.Range("J10").Value = 0.1683260544097508
$ws.Range("M10").Value = 0.7761303333333333
$ws.Range("N10").Value = 2.328391
$ws.Range("O10").Value = 0.03399380269149206
$ws.Range("P10").Value = 0.03399380269149207
$ws.Range("Q10").Value = 30.571024088098
$ws.Range("R10").Value = 275.139216792882
$ws.Range("S10").Value = 0.005722042681442426
$ws.Range("T10").Value = 0.005722042681442427
# Row 11
$ws.Range("I11").Value = 0.1683260544097508
$ws.Range("J11").Value = 0.1683260544097508
$ws.Range("M11").Value = 6.522593333333333
$ws.Range("N11").Value = 19.56778
$ws.Range("O11").Value = 0.2856836555503455
$ws.Range("P11").Value = 0.2856836555503455
$ws.Range("Q11").Value = 256.91865057484
$ws.Range("R11").Value = 2312.26785517356
$ws.Range("S11").Value = 0.04808800254814397
$ws.Range("T11").Value = 0.04808800254814397
# Row 12
$ws.Range("G12").Value = 44.61912266666666
$ws.Range("H12").Value = 133.857368
$ws.Range("I12").Value = 0.1906764423241422
$ws.Range("J12").Value = 0.1906764423241422
$ws.Range("M12").Value = 5.375839
$ws.Range("N12").Value = 16.127517
$ws.Range("O12").Value = 0.2354568587499626
$ws.Range("P12").Value = 0.2354568587499626
$ws.Range("Q12").Value = 239.8652197772506
$ws.Range("R12").Value = 2158.786977995256
$ws.Range("S12").Value = 0.04489607614726094
$ws.Range("T12").Value = 0.04489607614726094
# Row 13
$ws.Range("G13").Value = 44.61912266666666
$ws.Range("H13").Value = 133.857368
$ws.Range("I13").Value = 0.1906764423241422
$ws.Range("J13").Value = 0.1906764423241422
$ws.Range("O13").Value = 0.007131134316291014
$ws.Range("P13").Value = 0.007131134316291014
$ws.Range("Q13").Value = 7.264647583932444
$ws.Range("R13").Value = 65.381828255392
$ws.Range("S13").Value = 0.001359739321165975
$ws.Range("T13").Value = 0.001359739321165975
# Row 14
$ws.Range("G14").Value = 44.61912266666666
$ws.Range("H14").Value = 133.857368
$ws.Range("I14").Value = 0.1906764423241422
$ws.Range("J14").Value = 0.1906764423241422
$ws.Range("M14").Value = 9.994147
$ws.Range("N14").Value = 29.982441
$ws.Range("O14").Value = 0.4377345486919088
$ws.Range("P14").Value = 0.4377345486919088
$ws.Range("Q14").Value = 445.9300709416986
$ws.Range("R14").Value = 4013.370638475287
$ws.Range("S14").Value = 0.08346566642693717
$ws.Range("T14").Value = 0.08346566642693716
# Row 15
$ws.Range("G15").Value = 44.61912266666666
$ws.Range("H15").Value = 133.857368
$ws.Range("I15").Value = 0.1906764423241422
$ws.Range("J15").Value = 0.1906764423241422
$ws.Range("M15").Value = 0.7761303333333333
$ws.Range("N15").Value = 2.328391
$ws.Range("O15").Value = 0.03399380269149206
$ws.Range("P15").Value = 0.03399380269149207
$ws.Range("Q15").Value = 34.63025454832088
$ws.Range("R15").Value = 311.6722909348879
$ws.Range("S15").Value = 0.006481817358282557
$ws.Range("T15").Value = 0.006481817358282558
# Row 16
$ws.Range("G16").Value = 44.61912266666666
$ws.Range("H16").Value = 133.857368
$ws.Range("I16").Value = 0.1906764423241422
$ws.Range("J16").Value = 0.1906764423241422
$ws.Range("M16").Value = 6.522593333333333
$ws.Range("N16").Value = 19.56778
$ws.Range("O16").Value = 0.2856836555503455
$ws.Range("P16").Value = 0.2856836555503455
$ws.Range("Q16").Value = 291.0323920447822
$ws.Range("R16").Value = 2619.29152840304
$ws.Range("S16").Value = 0.05447314307049558
$ws.Range("T16").Value = 0.05447314307049557
# Row 17
$ws.Range("G17").Value = 18.49514433333333
$ws.Range("H17").Value = 55.485433
$ws.Range("I17").Value = 0.07903759892585487
$ws.Range("J17").Value = 0.07903759892585487
$ws.Range("M17").Value = 5.375839
$ws.Range("N17").Value = 16.127517
$ws.Range("O17").Value = 0.2354568587499626
$ws.Range("P17").Value = 0.2354568587499626
$ws.Range("Q17").Value = 99.42691821776233
$ws.Range("R17").Value = 894.842263959861
$ws.Range("S17").Value = 0.0186099447662212
$ws.Range("T17").Value = 0.0186099447662212
# Row 18
$ws.Range("G18").Value = 18.49514433333333
$ws.Range("H18").Value = 55.485433
$ws.Range("I18").Value = 0.07903759892585487
$ws.Range("J18").Value = 0.07903759892585487
$ws.Range("O18").Value = 0.007131134316291014
$ws.Range("P18").Value = 0.007131134316291014
$ws.Range("Q18").Value = 3.011280759583555
$ws.Range("R18").Value = 27.101526836252
$ws.Range("S18").Value = 0.0005636277339774095
$ws.Range("T18").Value = 0.0005636277339774095
# Row 19
$ws.Range("G19").Value = 18.49514433333333
$ws.Range("H19").Value = 55.485433
$ws.Range("I19").Value = 0.07903759892585487
$ws.Range("J19").Value = 0.07903759892585487
$ws.Range("M19").Value = 9.994147
$ws.Range("N19").Value = 29.982441
$ws.Range("O19").Value = 0.4377345486919088
$ws.Range("P19").Value = 0.4377345486919088
$ws.Range("Q19").Value = 184.8431912535503
$ws.Range("R19").Value = 1663.588721281953
$ws.Range("S19").Value = 0.03459748769550117
$ws.Range("T19").Value = 0.03459748769550117
# Row 20
$ws.Range("G20").Value = 18.49514433333333
$ws.Range("H20").Value = 55.485433
$ws.Range("I20").Value = 0.07903759892585487
$ws.Range("J20").Value = 0.07903759892585487
$ws.Range("M20").Value = 0.7761303333333333
$ws.Range("N20").Value = 2.328391
$ws.Range("O20").Value = 0.03399380269149206
$ws.Range("P20").Value = 0.03399380269149207
$ws.Range("Q20").Value = 14.35464253647811
$ws.Range("R20").Value = 129.191782828303
$ws.Range("S20").Value = 0.002686788543094795
$ws.Range("T20").Value = 0.002686788543094796
# Row 21
$ws.Range("G21").Value = 18.49514433333333
$ws.Range("H21").Value = 55.485433
$ws.Range("I21").Value = 0.07903759892585487
$ws.Range("J21").Value = 0.07903759892585487
$ws.Range("M21").Value = 6.522593333333333
$ws.Range("N21").Value = 19.56778
$ws.Range("O21").Value = 0.2856836555503455
$ws.Range("P21").Value = 0.2856836555503455
$ws.Range("Q21").Value = 120.6363051276378
$ws.Range("R21").Value = 1085.72674614874
$ws.Range("S21").Value = 0.02257975018706028
$ws.Range("T21").Value = 0.02257975018706028
# Row 22
$ws.Range("G22").Value = 52.59995866666666
$ws.Range("H22").Value = 157.799876
$ws.Range("I22").Value = 0.2247819406913095
$ws.Range("J22").Value = 0.2247819406913095
$ws.Range("M22").Value = 5.375839
$ws.Range("N22").Value = 16.127517
$ws.Range("O22").Value = 0.2354568587499626
$ws.Range("P22").Value = 0.2354568587499626
$ws.Range("Q22").Value = 282.7689091986546
$ws.Range("R22").Value = 2544.920182787892
$ws.Range("S22").Value = 0.05292644965889613
$ws.Range("T22").Value = 0.05292644965889613
# Row 23
$ws.Range("G23").Value = 52.59995866666666
$ws.Range("H23").Value = 157.799876
$ws.Range("I23").Value = 0.2247819406913095
$ws.Range("J23").Value = 0.2247819406913095
$ws.Range("O23").Value = 0.007131134316291014
$ws.Range("P23").Value = 0.007131134316291014
$ws.Range("Q23").Value = 8.564044736993777
$ws.Range("R23").Value = 77.076402632944
$ws.Range("S23").Value = 0.001602950210946289
$ws.Range("T23").Value = 0.001602950210946289
# Row 24
$ws.Range("G24").Value = 52.59995866666666
$ws.Range("H24").Value = 157.799876
$ws.Range("I24").Value = 0.2247819406913095
$ws.Range("J24").Value = 0.2247819406913095
$ws.Range("M24").Value = 9.994147
$ws.Range("N24").Value = 29.982441
$ws.Range("O24").Value = 0.4377345486919088
$ws.Range("P24").Value = 0.4377345486919088
$ws.Range("Q24").Value = 525.6917191085906
$ws.Range("R24").Value = 4731.225471977315
$ws.Range("S24").Value = 0.09839482136260178
$ws.Range("T24").Value = 0.09839482136260178
# Row 25
$ws.Range("G25").Value = 52.59995866666666
$ws.Range("H25").Value = 157.799876
$ws.Range("I25").Value = 0.2247819406913095
$ws.Range("J25").Value = 0.2247819406913095
$ws.Range("M25").Value = 0.7761303333333333
$ws.Range("N25").Value = 2.328391
$ws.Range("O25").Value = 0.03399380269149206
$ws.Range("P25").Value = 0.03399380269149207
$ws.Range("Q25").Value = 40.82442345327954
$ws.Range("R25").Value = 367.419811079516
$ws.Range("S25").Value = 0.007641192940471046
$ws.Range("T25").Value = 0.007641192940471048
# Row 26
$ws.Range("G26").Value = 52.59995866666666
$ws.Range("H26").Value = 157.799876
$ws.Range("I26").Value = 0.2247819406913095
$ws.Range("J26").Value = 0.2247819406913095
$ws.Range("M26").Value = 6.522593333333333
$ws.Range("N26").Value = 19.56778
$ws.Range("O26").Value = 0.2856836555503455
$ws.Range("P26").Value = 0.2856836555503455
$ws.Range("Q26").Value = 343.0881397328088
$ws.Range("R26").Value = 3087.79325759528
$ws.Range("S26").Value = 0.06421652651839427
$ws.Range("T26").Value = 0.06421652651839427
